# Update the worksheet with new data rows (302-328), extending data through 2021-07-25
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(44376,0,0,0),
    @(44377,0,0,0),
    @(44378,0,0,0),
    @(44379,0,0,0),
    @(44380,0,0,0),
    @(44381,0,0,0),
    @(44382,0,0,0),
    @(44383,0,0,0),
    @(44384,0,0,0),
    @(44385,0,0,0),
    @(44386,0,0,0),
    @(44387,0,0,0),
    @(44388,2,2,11.67201634082288),
    @(44389,0,2,11.67201634082288),
    @(44390,0,2,11.67201634082288),
    @(44391,0,2,11.67201634082288),
    @(44392,0,2,11.67201634082288),
    @(44393,1,3,17.50802451123432),
    @(44394,1,4,23.34403268164575),
    @(44395,0,2,11.67201634082288),
    @(44396,0,2,11.67201634082288),
    @(44397,0,2,11.67201634082288),
    @(44398,0,2,11.67201634082288),
    @(44399,1,3,17.50802451123432),
    @(44400,1,3,17.50802451123432),
    @(44401,1,3,17.50802451123432),
    @(44402,1,4,23.34403268164575)
)

$startRow = 302
$endRow = $startRow + $newRows.Count - 1

# Copy the date-column style/format from the last existing data row (A301) down
# to the new rows so formatting (bold, border, date number format, alignment) matches.
$ws.Range("A301").Copy()
$ws.Range("A" + $startRow + ":A" + $endRow).PasteSpecial(-4122)

$r = $startRow
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
